$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has the "description" header in row 2 (style 1) and
# the "format string" header in row 3 (style 2). Swap them (values AND
# formatting) so the format row comes first. A blank-cell-aware
# ClearContents before each paste ensures cells that should end up empty
# really do (Copy-into-destination otherwise leaves stale values behind
# in cells where the source range is blank).
$row2 = $ws.Range("B2:O2")
$row3 = $ws.Range("B3:O3")
$tempRow = $ws.Range("B60:O60")

$row2.Copy($tempRow)
$row2.ClearContents()
$row3.Copy($row2)
$row3.ClearContents()
$tempRow.Copy($row3)
$tempRow.Clear()

# Restore the selection to F20, as recorded in the saved workbook view.
[void]$ws.Range("F20").Select()
